# Updated symbol list on Tue Jan 17 20:35:29 UTC 2023 with GitHub Actions
# Applies refreshed Price (column D) and Volume(1h) (column E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.86"
$ws.Range("E2").Value = "'0.84%"
$ws.Range("D3").Value = "'32.16"
$ws.Range("E3").Value = "'1.30%"
$ws.Range("D4").Value = "'4.981"
$ws.Range("E4").Value = "'-3.10%"
$ws.Range("D5").Value = "'0.07910"
$ws.Range("E5").Value = "'-3.72%"
$ws.Range("D6").Value = "'2.108"
$ws.Range("E6").Value = "'-17.62%"
$ws.Range("D7").Value = "'7.860"
$ws.Range("E7").Value = "'0.04%"
$ws.Range("D8").Value = "'3.811"
$ws.Range("E8").Value = "'-1.23%"
$ws.Range("D9").Value = "'0.9269"
$ws.Range("E9").Value = "'-0.16%"
$ws.Range("D10").Value = "'0.1763"
$ws.Range("E10").Value = "'0.21%"
$ws.Range("D11").Value = "'0.08035"
$ws.Range("D12").Value = "'0.08775"
$ws.Range("E12").Value = "'-3.04%"
$ws.Range("E13").Value = "'4.80%"
$ws.Range("E14").Value = "'0.43%"
$ws.Range("D15").Value = "'0.001509"
$ws.Range("E15").Value = "'-1.02%"
$ws.Range("D16").Value = "'0.005997"
$ws.Range("E16").Value = "'0.05%"
$ws.Range("E17").Value = "'-4.13%"
$ws.Range("D18").Value = "'2.280"
$ws.Range("E18").Value = "'-0.22%"
$ws.Range("E19").Value = "'0.78%"
$ws.Range("D20").Value = "'0.1289"
$ws.Range("D21").Value = "'4.198"
$ws.Range("E21").Value = "'-1.28%"
$ws.Range("E22").Value = "'6.62%"
$ws.Range("D23").Value = "'0.04606"
$ws.Range("E23").Value = "'-0.53%"
$ws.Range("D24").Value = "'0.001236"
$ws.Range("E24").Value = "'-0.63%"
$ws.Range("D25").Value = "'0.004494"
$ws.Range("E25").Value = "'-1.48%"
$ws.Range("D26").Value = "'0.0001250"
$ws.Range("E26").Value = "'4.32%"
$ws.Range("D39").Value = "'0.01737"
$ws.Range("E39").Value = "'-2.61%"
$ws.Range("D40").Value = "'0.04802"
$ws.Range("E40").Value = "'4.08%"
$ws.Range("D41").Value = "'0.007346"
$ws.Range("E41").Value = "'6.59%"
$ws.Range("D42").Value = "'0.1368"
$ws.Range("E42").Value = "'-0.88%"
$ws.Range("D43").Value = "'0.002339"
$ws.Range("E43").Value = "'9.52%"
$ws.Range("D44").Value = "'0.01104"
$ws.Range("E44").Value = "'12.39%"
$ws.Range("D45").Value = "'0.00006015"
$ws.Range("E45").Value = "'-2.45%"
$ws.Range("E46").Value = "'0.31%"
$ws.Range("E47").Value = "'-59.52%"
$ws.Range("D48").Value = "'0.8204"
$ws.Range("E48").Value = "'2.09%"
$ws.Range("E49").Value = "'0.31%"
$ws.Range("E50").Value = "'0.31%"
